$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet2")

# ------------------------------------------------------------------
# 1) Replace the "O" (done marker) strings with "완료" in column A,
#    and make sure number-format/style matches the plain centered
#    style (same as A55/A58/A60 which already use that style).
# ------------------------------------------------------------------
$ws.Range("A55").Value = "완료"
$ws.Range("A58").Value = "완료"
$ws.Range("A60").Value = "완료"

$ws.Range("A56").Value = "완료"
$ws.Range("A55").Copy() | Out-Null
$ws.Range("A56").PasteSpecial(-4122) | Out-Null

$ws.Range("A57").Value = "완료"
$ws.Range("A55").Copy() | Out-Null
$ws.Range("A57").PasteSpecial(-4122) | Out-Null

$ws.Range("A62").Value = "완료"
$ws.Range("A55").Copy() | Out-Null
$ws.Range("A62").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Update the text of C63 (reuse - new data -> reuse - new data(bmi 데이터))
# ------------------------------------------------------------------
$ws.Range("C63").Value = "reuse - new data(bmi 데이터)"

# ------------------------------------------------------------------
# 3) Row 64 date moved from 43075 to 43085
# ------------------------------------------------------------------
$ws.Range("A64").Value = 43085

# ------------------------------------------------------------------
# 4) Row 65 date moved from 43075 to 43084, and gets a new owner
#    "이은" in column B (new cell), plus the highlight below.
# ------------------------------------------------------------------
$ws.Range("A65").Value = 43084
$ws.Range("B65").Value = "이은"

# ------------------------------------------------------------------
# 5) Row 63 gets a new owner "영택" in column B (new cell).
# ------------------------------------------------------------------
$ws.Range("B63").Value = "영택"

# ------------------------------------------------------------------
# 6) Highlight (yellow fill) rows 59, 63, 65, 68, 69 in columns A:C
#    (this is the new "done recently" highlight introduced by the edit)
# ------------------------------------------------------------------
$yellow = 65535
foreach ($r in 59,63,65,68,69) {
    $ws.Range("A" + $r).Interior.Color = $yellow
    $ws.Range("B" + $r).Interior.Color = $yellow
    $ws.Range("C" + $r).Interior.Color = $yellow
}

# ------------------------------------------------------------------
# 7) Update the active selection shown in the sheet view
# ------------------------------------------------------------------
$ws.Range("B64").Select()

